# "Nueva tarea - reportes de devolucion"
#
# Adds a new task row to the bottom of the Hoja1 task list, records a
# 50% progress value for the "comprar CPU" task (row 11), and updates
# the sheet's view/selection to reflect where the user ended up after
# making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "comprar CPU" (row 11) now has a 50% progress value in column C,
# formatted the same way as the other percentage cells in the sheet.
$ws.Range("C11").Value = 0.5
$ws.Range("C11").NumberFormat = "0%"

# New task appended at the end of the list (row 43).
$ws.Range("A43").Value = "revisar reportes de venta - devolucion"

# Scroll the view down and leave the selection on the row just past the
# new task, matching where the user was working.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A44").Select()
